$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.551.60"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "1.881.59"
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "243.70"
$ws.Range("E5").Value = "  +4.75%  "

$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("D8").Value = "42.99"
$ws.Range("E8").Value = "  +5.31%  "

$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("D11").Value = "0.0993"
$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").Value = "2.153.14"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("E13").Value = "  +7.93%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.686"
$ws.Range("E14").Value = "  +1.99%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.846.67"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").Value = "4.78"
$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Value = "35.555.92"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "71.51"
$ws.Range("E18").Value = "  +2.18%  "

$ws.Range("E19").Value = "  +2.22%  "

$ws.Range("D20").Value = "243.68"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("E22").Value = "  +2.40%  "

$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").Value = "171.19"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("E26").Value = "  +26.16%  "

$ws.Range("D27").Value = "8.35"
$ws.Range("E27").Value = "  +6.46%  "

$ws.Range("D28").Value = "17.88"
$ws.Range("E28").Value = "  +2.35%  "

$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").Value = "0.950"
$ws.Range("E30").Value = "  +27.10%  "

$ws.Range("D31").Value = "0.0564"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("E32").Value = "  +3.02%  "

$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("E34").Value = "  +3.42%  "

$ws.Range("E35").Value = "  +11.57%  "

$ws.Range("D36").Value = "2.06"
$ws.Range("E36").Value = "  +5.52%  "

$ws.Range("D37").Value = "1.35"
$ws.Range("E37").Value = "  +10.97%  "

$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").Value = "90.41"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").Value = "1.357.85"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "15.31"
$ws.Range("E42").Value = "  +3.93%  "

$ws.Range("E43").Value = "  +11.83%  "

$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "12.99"
$ws.Range("E44").Value = "  +44.38%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "2.34"
$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "6.71"
$ws.Range("E47").Value = "  +6.29%  "

$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "45.63"
$ws.Range("E48").Value = "  +34.51%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "2.072.88"
$ws.Range("E50").Value = "  +2.08%  "

$ws.Range("E51").Value = "  +2.61%  "
